$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W8)
$ws1.Range("D2").Value = 640
$ws1.Range("H2").Value = 9.460000000000001
$ws1.Range("L2").Value = 0.82

# Row 3 (W9)
$ws1.Range("D3").Value = 632
$ws1.Range("H3").Value = 8.57
$ws1.Range("L3").Value = 0.98

# Row 4 (W10)
$ws1.Range("D4").Value = 626
$ws1.Range("H4").Value = 7.64
$ws1.Range("L4").Value = 0.83

# Row 5 (W11)
$ws1.Range("D5").Value = 621
$ws1.Range("H5").Value = 6.7

# Row 6 (W12)
$ws1.Range("D6").Value = 611
$ws1.Range("H6").Value = 5.79
$ws1.Range("L6").Value = 1.16

# Row 7 (W13)
$ws1.Range("D7").Value = 596
$ws1.Range("H7").Value = 4.91
$ws1.Range("L7").Value = 0.85

# Row 8 (W14)
$ws1.Range("D8").Value = 596
$ws1.Range("H8").Value = 3.91
$ws1.Range("L8").Value = 1.12

# Row 9 (W15)
$ws1.Range("D9").Value = 590
$ws1.Range("H9").Value = 2.94
$ws1.Range("L9").Value = 0.97

# Row 10 (W16)
$ws1.Range("D10").Value = 582
$ws1.Range("H10").Value = 1.97
$ws1.Range("L10").Value = 1.03

# Row 11 (W17)
$ws1.Range("D11").Value = 581
$ws1.Range("H11").Value = 0.97
$ws1.Range("I11").Value = "Low"
$ws1.Range("L11").Value = 1.17

# Row 12 (W18)
$ws1.Range("D12").Value = 575
$ws1.Range("L12").Value = 0.86

# Row 13 (W19)
$ws1.Range("D13").Value = 575
$ws1.Range("L13").Value = 1.19

# Row 14 (W20)
$ws1.Range("D14").Value = 569

# Row 15 (W21)
$ws1.Range("D15").Value = 560
$ws1.Range("L15").Value = 1.05

# Row 16 (W22)
$ws1.Range("D16").Value = 557
$ws1.Range("L16").Value = 1.09

# Row 17 (W23)
$ws1.Range("D17").Value = 546
$ws1.Range("L17").Value = 1.17

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")

# These values are stored as text (not numbers) in the sheet, so force the
# cell format to Text before assigning, otherwise Excel would auto-convert
# the numeric-looking string into a real number.
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "9457"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "4912"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "2519"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "640"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "546"
